$wb = $excel.ActiveWorkbook

# Source sheet already contains the "checked" (P / Wingdings2 checkmark) style
# at these exact coordinates - reuse its formatting (style 18) via a
# formats-only paste, then stamp the "P" value (shared with other "checked"
# cells elsewhere in the workbook).
$srcWs = $wb.Worksheets.Item("Hoja1 (2)")
$dstWs = $wb.Worksheets.Item("Hoja1 (4)")

# destination address -> source address (same cell in the template sheet,
# except C5 which is blank there; D5 carries the identical "checked" style).
$targets = @{
    "C4"  = "C4"
    "E4"  = "E4"
    "F4"  = "F4"
    "C5"  = "D5"
    "D5"  = "D5"
    "D6"  = "D6"
    "E8"  = "E8"
    "F8"  = "F8"
    "D9"  = "D9"
    "E9"  = "E9"
    "D11" = "D11"
    "E11" = "E11"
}

foreach ($addr in $targets.Keys) {
    $srcAddr = $targets[$addr]
    $src = $srcWs.Range($srcAddr)
    $dst = $dstWs.Range($addr)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = "P"
}

$excel.CutCopyMode = 0

# Restore the selections left behind by the editing session, per sheet.
$wb.Worksheets.Item("Hoja1 (2)").Range("C4").Select()
$wb.Worksheets.Item("Hoja1 (3)").Range("E9").Select()
$dstWs.Activate()
$dstWs.Range("E18").Select()
